$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143 (shifts existing rows 143-146 down to 144-147)
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with the new record.
# Columns A, B, C, E, F, G, H, I, J, K, L, R are the same constant values
# used throughout this block of rows (Vega Modelo de Temuco / Papaya entries).
$ws.Range("A143").Value = 10
$ws.Range("B143").Value = "Vega Modelo de Temuco"
$ws.Range("C143").Value = "La Araucanía"
$ws.Range("D143").Value = 45265
$ws.Range("E143").Value = 9
$ws.Range("F143").Value = "Fruta"
$ws.Range("G143").Value = 100108
$ws.Range("H143").Value = "Tropicales y subtropicales"
$ws.Range("I143").Value = 100108004
$ws.Range("J143").Value = "Papaya"
$ws.Range("K143").Value = "Cultivar IV Región"
$ws.Range("L143").Value = "Primera"
$ws.Range("M143").Value = 120
$ws.Range("N143").Value = 37000
$ws.Range("O143").Value = 37000
$ws.Range("P143").Value = 37000
$ws.Range("Q143").Value = "$/caja 15 kilos granel"
$ws.Range("R143").Value = "Provincia del Elquí"
$ws.Range("S143").Value = 2467
$ws.Range("T143").Value = 15
